# Actualiza base de datos EC: rota los registros de los 3 trabajadores
# en las filas 16-18 de la hoja "Hoja1" (cada fila pasa a mostrar los
# datos que antes tenia la fila siguiente, de forma ciclica):
#   nueva fila 16 = vieja fila 17 (NILSON JESUS GARCIA YEPES)
#   nueva fila 17 = vieja fila 18 (EDILBERTO QUINTERO CUELLAR)
#   nueva fila 18 = vieja fila 16 (MARIA ALEJANDRA MACIAS QUINTERO)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

# Lee los valores actuales (antes de sobreescribir nada)
$c16 = $ws.Range("C16").Value2
$d16 = $ws.Range("D16").Value2
$e16 = $ws.Range("E16").Value2
$f16 = $ws.Range("F16").Value2
$g16 = $ws.Range("G16").Value2

$c17 = $ws.Range("C17").Value2
$d17 = $ws.Range("D17").Value2
$e17 = $ws.Range("E17").Value2
$f17 = $ws.Range("F17").Value2
$g17 = $ws.Range("G17").Value2

$c18 = $ws.Range("C18").Value2
$d18 = $ws.Range("D18").Value2
$e18 = $ws.Range("E18").Value2
$f18 = $ws.Range("F18").Value2
$g18 = $ws.Range("G18").Value2

# Fila 16 <- datos de la vieja fila 17
$ws.Range("C16").Value2 = $c17
$ws.Range("D16").Value2 = $d17
$ws.Range("E16").Value2 = $e17
$ws.Range("F16").Value2 = $f17
$ws.Range("G16").Value2 = $g17

# Fila 17 <- datos de la vieja fila 18
$ws.Range("C17").Value2 = $c18
$ws.Range("D17").Value2 = $d18
$ws.Range("E17").Value2 = $e18
$ws.Range("F17").Value2 = $f18
$ws.Range("G17").Value2 = $g18

# Fila 18 <- datos de la vieja fila 16
$ws.Range("C18").Value2 = $c16
$ws.Range("D18").Value2 = $d16
$ws.Range("E18").Value2 = $e16
$ws.Range("F18").Value2 = $f16
$ws.Range("G18").Value2 = $g16
